$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("0.06059416010975838", "0.9823344945907593", "0.02161066979169846", "0.9961565136909485"),
    @("0.007842882536351681", "0.9984692335128784", "0.01125249546021223", "0.9961565136909485"),
    @("0.003263595979660749", "0.9990850687026978", "0.003971421625465155", "0.99895179271698"),
    @("0.002217731904238462", "0.9995073080062866", "0.001151856267824769", "0.99895179271698"),
    @("0.001602303236722946", "0.9995601177215576", "0.0007560593658126891", "1"),
    @("0.001637882320210338", "0.9995601177215576", "0.001006737584248185", "0.9996505975723267"),
    @("0.000885470537468791", "0.9997184872627258", "0.0005875350907444954", "0.9996505975723267"),
    @("0.0007942115771584213", "0.9997360706329346", "6.64602848701179E-05", "1"),
    @("0.0007708691409789026", "0.9998416304588318", "4.504147364059463E-05", "1"),
    @("0.0007734073442406952", "0.9998064637184143", "0.001351145794615149", "0.9996505975723267"),
    @("0.0005454458296298981", "0.999876856803894", "4.096049451618455E-05", "1"),
    @("0.0004644978325814009", "0.9998944401741028", "0.0001232130452990532", "1"),
    @("0.0003671460435725749", "0.999876856803894", "3.892044333042577E-05", "1"),
    @("0.0004789872036781162", "0.999876856803894", "3.573861249606125E-05", "1"),
    @("0.000153091357788071", "0.9999296069145203", "3.522664337651804E-05", "1"),
    @("0.0006993327406235039", "0.9998064637184143", "0.00133404484950006", "0.9996505975723267"),
    @("2.8134945750935E-05", "1", "0.003296397626399994", "0.9996505975723267"),
    @("0.0008119892445392907", "0.9998064637184143", "1.162749504146632E-05", "1"),
    @("0.0001912672450998798", "0.9999824166297913", "6.207340629771352E-05", "1"),
    @("0.0002940595149993896", "0.9998944401741028", "1.420518674422055E-05", "1"),
    @("0.0005880878889001906", "0.9998064637184143", "3.545456638676114E-05", "1"),
    @("0.0001963468384929001", "0.9999824166297913", "0.001376168569549918", "0.9996505975723267"),
    @("0.000145534664625302", "0.9999648332595825", "0.0005252000992186368", "0.9996505975723267"),
    @("0.0004788801015820354", "0.999876856803894", "0.0001972745230887085", "1"),
    @("0.0001336545683443546", "0.9999296069145203", "0.0001030836938298307", "1"),
    @("0.000238155378610827", "0.9999296069145203", "7.781247404636815E-06", "1"),
    @("0.0002474871871527284", "0.9999648332595825", "2.518243491067551E-06", "1"),
    @("9.967148071154952E-05", "0.9999648332595825", "1.216146330307311E-07", "1"),
    @("3.691367282954161E-06", "1", "1.957227908633286E-07", "1"),
    @("0.0002706167288124561", "0.999947190284729", "7.033373776721419E-07", "1"),
    @("0.0003828966291621327", "0.9998944401741028", "0.001436293474398553", "0.9996505975723267"),
    @("0.0001377027074340731", "0.999947190284729", "1.339682512480067E-06", "1"),
    @("4.263108712621033E-05", "0.9999648332595825", "1.15747802453825E-07", "1"),
    @("0.0002982307341881096", "0.9999648332595825", "7.663884815656274E-08", "1"),
    @("2.516223867132794E-05", "0.9999824166297913", "1.307878072509538E-08", "1"),
    @("0.0004315355618018657", "0.999876856803894", "1.40689550676143E-07", "1"),
    @("0.0001247839245479554", "0.9999648332595825", "6.688231337648176E-07", "1"),
    @("0.0002747255493886769", "0.9999648332595825", "5.793778967699836E-08", "1"),
    @("0.0002614011755213141", "0.9999120235443115", "2.244458983113873E-07", "1"),
    @("0.0002011915930779651", "0.999947190284729", "1.091835565603105E-05", "1"),
    @("0.0002317297767149284", "0.999947190284729", "1.025874141191707E-07", "1"),
    @("0.0006185008678585291", "0.9999824166297913", "2.322690988876275E-06", "1"),
    @("0.0006156565505079925", "0.9998592138290405", "1.612890514479659E-07", "1"),
    @("0.0003933866391889751", "0.9999296069145203", "3.656779881566763E-06", "1"),
    @("0.0001864439836936072", "0.999947190284729", "1.475894322311433E-07", "1"),
    @("0.0002747914113570005", "0.9999648332595825", "5.239343863649992E-06", "1"),
    @("0.0001239181874552742", "0.9999648332595825", "3.057483297652652E-07", "1"),
    @("0.000210453275940381", "0.9999648332595825", "3.061303459617193E-06", "1"),
    @("0.0007493131561204791", "0.999876856803894", "2.777996428449114E-07", "1"),
    @("0.0001434707082808018", "0.9999648332595825", "5.985158679777669E-08", "1"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = [double]$data[$i][0]
    $ws.Cells.Item($row, 2).Value = [double]$data[$i][1]
    $ws.Cells.Item($row, 3).Value = [double]$data[$i][2]
    $ws.Cells.Item($row, 4).Value = [double]$data[$i][3]
}
